$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 21

$ws.Cells.Item($row, 1).Value = "2025-08-16 13:02:40 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-16 18:32:40 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""

$ws.Range("A21:H21").HorizontalAlignment = -4108
$ws.Range("A21:H21").VerticalAlignment = -4108
